# Updates the crypto price/volume table (columns D and E) with refreshed
# figures from the GitHub Actions scrape, matching the target commit.
# Rows 48/49 also swap coin identity (RenderToken <-> BabyDogeCoin).
#
# Price cells (column D) often look numeric (e.g. "0.9997", "243.20") but
# must stay plain text, exactly as in the source data - so NumberFormat is
# forced to "@" (Text) right before assigning those values to stop Excel
# from silently re-interpreting them as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.322.83'
$ws.Range("E2").Value = '  +1.79%  '

$ws.Range("D3").Value = '1.841.10'
$ws.Range("E3").Value = '  +0.67%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9989'
$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.20'
$ws.Range("E5").Value = '  -0.64%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6868'
$ws.Range("E6").Value = '  -0.41%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9997'
$ws.Range("E7").Value = '  -0.07%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3032'
$ws.Range("E8").Value = '  -0.14%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07531'
$ws.Range("E9").Value = '  -1.46%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.29'
$ws.Range("E10").Value = '  +0.32%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07654'
$ws.Range("E11").Value = '  -1.66%  '

$ws.Range("D12").Value = '1.842.42'
$ws.Range("E12").Value = '  +0.70%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.088'
$ws.Range("E13").Value = '  +0.26%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6875'
$ws.Range("E14").Value = '  +1.54%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '89.63'
$ws.Range("E15").Value = '  -3.71%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.307'
$ws.Range("E16").Value = '  -2.04%  '

$ws.Range("D17").Value = '29.313.46'
$ws.Range("E17").Value = '  +1.79%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008261'
$ws.Range("E18").Value = '  +0.55%  '

$ws.Range("E19").Value = '  +0.86%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '234.92'
$ws.Range("E20").Value = '  -2.53%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.60'
$ws.Range("E21").Value = '  -0.32%  '

$ws.Range("E22").Value = '  -0.04%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.490'
$ws.Range("E23").Value = '  +0.82%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.9991'
$ws.Range("E24").Value = '  -0.07%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1457'
$ws.Range("E25").Value = '  -1.96%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '160.10'
$ws.Range("E26").Value = '  -1.12%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.854'
$ws.Range("E27").Value = '  +1.58%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.13'
$ws.Range("E28").Value = '  -0.30%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.521'
$ws.Range("E29").Value = '  -1.03%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.235'
$ws.Range("E30").Value = '  +0.44%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.140'
$ws.Range("E31").Value = '  -0.29%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.202'
$ws.Range("E32").Value = '  +1.21%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05140'
$ws.Range("E33").Value = '  +0.55%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7717'
$ws.Range("E34").Value = '  +0.13%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.850'
$ws.Range("E35").Value = '  +0.40%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.140'
$ws.Range("E36").Value = '  +0.47%  '

$ws.Range("D38").Value = '1.296.75'
$ws.Range("E38").Value = '  +3.31%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01847'
$ws.Range("E39").Value = '  -0.23%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.700'
$ws.Range("E40").Value = '  +0.16%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9442'
$ws.Range("E41").Value = '  -1.35%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '105.67'
$ws.Range("E42").Value = '  -0.93%  '

$ws.Range("E43").Value = '  -0.06%  '

$ws.Range("E44").Value = '  -6.17%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '9.690'
$ws.Range("E45").Value = '  +0.49%  '

$ws.Range("D46").Value = '1.989.96'
$ws.Range("E46").Value = '  +0.84%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5205'
$ws.Range("E47").Value = '  +1.04%  '

$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.771'
$ws.Range("E48").Value = '  +1.67%  '

$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00000000120'
$ws.Range("E49").Value = '  +0.00%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '63.31'
$ws.Range("E50").Value = '  -0.76%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05925'
$ws.Range("E51").Value = '  +0.85%  '
